# Auto-generated edit script: applies numeric cell updates to the
# Titan_Profits workbook sheets (ALC, ARM, CRP, CUL, GSM, LTW, WVR).
# Values/cells below mirror the authoritative OOXML diff exactly,
# including the few rows where a trailing column cell is removed
# (e.g. N132 on ALC) or newly introduced (e.g. N92 on CRP).
$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 163.66667
$ws.Range("I12").Value = 95.5
$ws.Range("J12").Value = 300
$ws.Range("K12").Value = 95.5
$ws.Range("L12").Value = 300
$ws.Range("M12").Value = 74.5
$ws.Range("N12").Value = -640
$ws.Range("H15").Value = 100094.91
$ws.Range("I15").Value = 100094.91
$ws.Range("K15").Value = 300284.73
$ws.Range("M15").Value = -300115.73
$ws.Range("H17").Value = 2404.1455
$ws.Range("J17").Value = 2404.1455
$ws.Range("L17").Value = 7212.4365
$ws.Range("N17").Value = -7548.4365
$ws.Range("H98").Value = 1126002.5
$ws.Range("I98").Value = 1250201
$ws.Range("K98").Value = 1250201
$ws.Range("M98").Value = -1248703
$ws.Range("H122").Value = 1126002.5
$ws.Range("I122").Value = 1250201
$ws.Range("K122").Value = 3750603
$ws.Range("M122").Value = -3748153
$ws.Range("H132").Value = 23420.043
$ws.Range("I132").Value = 23420.043
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 70260.129
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -67730.129
$ws.Range("N132").ClearContents()
$ws.Range("H135").Value = 1432
$ws.Range("I135").Value = 1664.5
$ws.Range("K135").Value = 14980.5
$ws.Range("M135").Value = -12445.5
$ws.Range("H137").Value = 23257256
$ws.Range("I137").Value = 34483840
$ws.Range("J137").Value = 2193.0715
$ws.Range("K137").Value = 103451520
$ws.Range("L137").Value = 6579.2145
$ws.Range("M137").Value = -103448970
$ws.Range("N137").Value = -11679.2145

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H4").Value = 1107.9166
$ws.Range("I4").Value = 1449.5
$ws.Range("K4").Value = 1449.5
$ws.Range("M4").Value = -1333.5
$ws.Range("H32").Value = 19449.293
$ws.Range("I32").Value = 4576.0146
$ws.Range("J32").Value = 163932.58
$ws.Range("K32").Value = 4576.0146
$ws.Range("L32").Value = 163932.58
$ws.Range("M32").Value = -4289.0146
$ws.Range("N32").Value = -164506.58
$ws.Range("H110").Value = 1301.6
$ws.Range("I110").Value = 1429
$ws.Range("J110").Value = 1004.3333
$ws.Range("K110").Value = 1429
$ws.Range("L110").Value = 1004.3333
$ws.Range("M110").Value = 616
$ws.Range("N110").Value = -5094.3333
$ws.Range("H132").Value = 2914.5
$ws.Range("I132").Value = 2493.6775
$ws.Range("K132").Value = 7481.032499999999
$ws.Range("M132").Value = -4951.032499999999

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1597.4286
$ws.Range("I16").Value = 1369
$ws.Range("J16").Value = 2008.6
$ws.Range("K16").Value = 1369
$ws.Range("L16").Value = 2008.6
$ws.Range("M16").Value = -1082
$ws.Range("N16").Value = -2582.6
$ws.Range("H31").Value = 3169.0356
$ws.Range("I31").Value = 1391.6216
$ws.Range("K31").Value = 1391.6216
$ws.Range("M31").Value = -1096.6216
$ws.Range("H34").Value = 3169.0356
$ws.Range("I34").Value = 1391.6216
$ws.Range("K34").Value = 1391.6216
$ws.Range("M34").Value = -1189.6216
$ws.Range("H58").Value = 1744.1515
$ws.Range("I58").Value = 1136.6
$ws.Range("J58").Value = 3642.75
$ws.Range("K58").Value = 1136.6
$ws.Range("L58").Value = 3642.75
$ws.Range("M58").Value = -933.5999999999999
$ws.Range("N58").Value = -4048.75
$ws.Range("H92").Value = 30000
$ws.Range("J92").Value = 30000
$ws.Range("L92").Value = 30000
$ws.Range("N92").Value = -34992
$ws.Range("H105").Value = 914.4545000000001
$ws.Range("I105").Value = 907.6667
$ws.Range("J105").Value = 945
$ws.Range("K105").Value = 907.6667
$ws.Range("L105").Value = 945
$ws.Range("M105").Value = 839.3333
$ws.Range("N105").Value = -4439
$ws.Range("H113").Value = 1597.4286
$ws.Range("I113").Value = 1369
$ws.Range("J113").Value = 2008.6
$ws.Range("K113").Value = 1369
$ws.Range("L113").Value = 2008.6
$ws.Range("M113").Value = 801
$ws.Range("N113").Value = -6348.6
$ws.Range("H122").Value = 1917.1786
$ws.Range("I122").Value = 1167.3
$ws.Range("J122").Value = 3791.875
$ws.Range("K122").Value = 3501.9
$ws.Range("L122").Value = 11375.625
$ws.Range("M122").Value = -1051.9
$ws.Range("N122").Value = -16275.625
$ws.Range("H136").Value = 1744.1515
$ws.Range("I136").Value = 1136.6
$ws.Range("J136").Value = 3642.75
$ws.Range("K136").Value = 3409.8
$ws.Range("L136").Value = 10928.25
$ws.Range("M136").Value = -859.7999999999997
$ws.Range("N136").Value = -16028.25

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H98").Value = 500
$ws.Range("I98").Value = 0
$ws.Range("J98").Value = 500
$ws.Range("K98").Value = 0
$ws.Range("L98").Value = 1500
$ws.Range("N98").Value = -4496
$ws.Range("M98").ClearContents()
$ws.Range("H131").Value = 10418557
$ws.Range("I131").Value = 633.3333
$ws.Range("J131").Value = 11496274
$ws.Range("K131").Value = 1899.9999
$ws.Range("L131").Value = 34488822
$ws.Range("M131").Value = 3140.0001
$ws.Range("N131").Value = -34498902
$ws.Range("H136").Value = 2860.0637
$ws.Range("I136").Value = 2205.6
$ws.Range("K136").Value = 6616.799999999999
$ws.Range("M136").Value = -1516.799999999999
$ws.Range("H138").Value = 1402.5
$ws.Range("I138").Value = 866.6667
$ws.Range("J138").Value = 3010
$ws.Range("K138").Value = 2600.0001
$ws.Range("L138").Value = 9030
$ws.Range("M138").Value = 2539.9999
$ws.Range("N138").Value = -19310

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H113").Value = 1850
$ws.Range("I113").Value = 1850
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 1850
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = 320
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 3558.0322
$ws.Range("I132").Value = 3077.1365
$ws.Range("J132").Value = 4733.5557
$ws.Range("K132").Value = 9231.4095
$ws.Range("L132").Value = 14200.6671
$ws.Range("M132").Value = -6701.4095
$ws.Range("N132").Value = -19260.6671
$ws.Range("H136").Value = 27150
$ws.Range("J136").Value = 27854
$ws.Range("L136").Value = 83562
$ws.Range("N136").Value = -88662

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 14900.429
$ws.Range("I22").Value = 475.25
$ws.Range("J22").Value = 34134
$ws.Range("K22").Value = 475.25
$ws.Range("L22").Value = 34134
$ws.Range("M22").Value = -180.25
$ws.Range("N22").Value = -34724
$ws.Range("H27").Value = 14900.429
$ws.Range("I27").Value = 475.25
$ws.Range("J27").Value = 34134
$ws.Range("K27").Value = 475.25
$ws.Range("L27").Value = 34134
$ws.Range("M27").Value = -368.25
$ws.Range("N27").Value = -34348
$ws.Range("H40").Value = 3048.1924
$ws.Range("I40").Value = 2206.625
$ws.Range("J40").Value = 3422.2222
$ws.Range("K40").Value = 2206.625
$ws.Range("L40").Value = 3422.2222
$ws.Range("M40").Value = -2070.625
$ws.Range("N40").Value = -3694.2222
$ws.Range("H68").Value = 1714.2858
$ws.Range("I68").Value = 1550
$ws.Range("K68").Value = 1550
$ws.Range("M68").Value = -801
$ws.Range("H71").Value = 1714.2858
$ws.Range("I71").Value = 1550
$ws.Range("K71").Value = 7750
$ws.Range("M71").Value = -4006
$ws.Range("H100").Value = 2560.6572
$ws.Range("I100").Value = 1892.1111
$ws.Range("J100").Value = 2792.077
$ws.Range("K100").Value = 1892.1111
$ws.Range("L100").Value = 2792.077
$ws.Range("M100").Value = -1351.1111
$ws.Range("N100").Value = -3874.077
$ws.Range("H122").Value = 3299.6843
$ws.Range("I122").Value = 2612.5334
$ws.Range("J122").Value = 3747.8262
$ws.Range("K122").Value = 7837.600199999999
$ws.Range("L122").Value = 11243.4786
$ws.Range("M122").Value = -5387.600199999999
$ws.Range("N122").Value = -16143.4786

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 3377.2285
$ws.Range("I81").Value = 2261
$ws.Range("K81").Value = 4522
$ws.Range("M81").Value = -3461
$ws.Range("H84").Value = 3377.2285
$ws.Range("I84").Value = 2261
$ws.Range("K84").Value = 22610
$ws.Range("M84").Value = -17306
$ws.Range("H107").Value = 887.9167
$ws.Range("I107").Value = 975.2
$ws.Range("K107").Value = 2925.6
$ws.Range("M107").Value = -1005.6

